# Update column F (dSF) values for the specified rows to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 4
    5  = -1
    7  = 5
    8  = 1
    9  = 4
    10 = -1
    11 = -1
    12 = -5
    13 = 2
    15 = -2
    16 = 2
    17 = 7
    18 = 2
    20 = 6
    22 = -4
    23 = -3
    24 = -1
    26 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
